# Applies the scheduled-runner price refresh to the Leve profit tables.
# Values come from an updated market-price snapshot; columns are plain
# data (no formulas) so each touched cell is written directly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value2 = 3619.261
$ws.Range("I28").Value2 = 5206.8335
$ws.Range("K28").Value2 = 5206.8335
$ws.Range("M28").Value2 = -4721.8335
$ws.Range("H62").Value2 = 4702.615
$ws.Range("I62").Value2 = 4694.9165
$ws.Range("J62").Value2 = 4795
$ws.Range("K62").Value2 = 4694.9165
$ws.Range("L62").Value2 = 4795
$ws.Range("M62").Value2 = -4070.9165
$ws.Range("N62").Value2 = -6043
$ws.Range("H65").Value2 = 4702.615
$ws.Range("I65").Value2 = 4694.9165
$ws.Range("J65").Value2 = 4795
$ws.Range("K65").Value2 = 23474.5825
$ws.Range("L65").Value2 = 23975
$ws.Range("M65").Value2 = -20354.5825
$ws.Range("N65").Value2 = -30215
$ws.Range("H70").Value2 = 22223954
$ws.Range("J70").Value2 = 1925.75
$ws.Range("L70").Value2 = 5777.25
$ws.Range("N70").Value2 = -6317.25
$ws.Range("H73").Value2 = 22223954
$ws.Range("J73").Value2 = 1925.75
$ws.Range("L73").Value2 = 5777.25
$ws.Range("N73").Value2 = -7649.25
$ws.Range("H92").Value2 = 577.5217
$ws.Range("I92").Value2 = 621.7
$ws.Range("K92").Value2 = 621.7
$ws.Range("M92").Value2 = 626.3
$ws.Range("H107").Value2 = 9373.286
$ws.Range("I107").Value2 = 9578.923000000001
$ws.Range("J107").Value2 = 6700
$ws.Range("K107").Value2 = 9578.923000000001
$ws.Range("L107").Value2 = 6700
$ws.Range("M107").Value2 = -7658.923000000001
$ws.Range("N107").Value2 = -10540
$ws.Range("H113").Value2 = 24865.8
$ws.Range("J113").Value2 = 38000
$ws.Range("L113").Value2 = 38000
$ws.Range("N113").Value2 = -44508
$ws.Range("H116").Value2 = 12225592
$ws.Range("J116").Value2 = 4283.3335
$ws.Range("L116").Value2 = 4283.3335
$ws.Range("N116").Value2 = -11167.3335
$ws.Range("H132").Value2 = 4719.6904
$ws.Range("I132").Value2 = 4740.8613
$ws.Range("K132").Value2 = 14222.5839
$ws.Range("M132").Value2 = -11692.5839

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value2 = 2999
$ws.Range("J23").Value2 = 2999
$ws.Range("L23").Value2 = 2999
$ws.Range("N23").Value2 = -3517
$ws.Range("H32").Value2 = 5891.5137
$ws.Range("I32").Value2 = 5891.5137
$ws.Range("K32").Value2 = 5891.5137
$ws.Range("M32").Value2 = -5604.5137
$ws.Range("H102").Value2 = 13746.111
$ws.Range("I102").Value2 = 20549.908
$ws.Range("J102").Value2 = 3054.4285
$ws.Range("K102").Value2 = 20549.908
$ws.Range("L102").Value2 = 3054.4285
$ws.Range("M102").Value2 = -18927.908
$ws.Range("N102").Value2 = -6298.4285
$ws.Range("H110").Value2 = 2330.2727
$ws.Range("I110").Value2 = 988.8333
$ws.Range("K110").Value2 = 988.8333
$ws.Range("M110").Value2 = 1056.1667
$ws.Range("H121").Value2 = 70955
$ws.Range("J121").Value2 = 70955
$ws.Range("L121").Value2 = 70955
$ws.Range("N121").Value2 = -74449
$ws.Range("H132").Value2 = 5186.84
$ws.Range("I132").Value2 = 5055.0264
$ws.Range("J132").Value2 = 5604.25
$ws.Range("K132").Value2 = 15165.0792
$ws.Range("L132").Value2 = 16812.75
$ws.Range("M132").Value2 = -12635.0792
$ws.Range("N132").Value2 = -21872.75
$ws.Range("H135").Value2 = 0
$ws.Range("J135").Value2 = 0
$ws.Range("L135").Value2 = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 2559.9524
$ws.Range("I20").Value2 = 1677.5
$ws.Range("J20").Value2 = 4324.857
$ws.Range("K20").Value2 = 1677.5
$ws.Range("L20").Value2 = 4324.857
$ws.Range("M20").Value2 = -1430.5
$ws.Range("N20").Value2 = -4818.857
$ws.Range("H99").Value2 = 18519.709
$ws.Range("I99").Value2 = 29370.54
$ws.Range("J99").Value2 = 5696
$ws.Range("K99").Value2 = 29370.54
$ws.Range("L99").Value2 = 5696
$ws.Range("M99").Value2 = -27872.54
$ws.Range("N99").Value2 = -8692
$ws.Range("H105").Value2 = 59030.723
$ws.Range("I105").Value2 = 78777.62
$ws.Range("K105").Value2 = 78777.62
$ws.Range("M105").Value2 = -77030.62

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1995.2858
$ws.Range("I16").Value2 = 1994.8334
$ws.Range("K16").Value2 = 1994.8334
$ws.Range("M16").Value2 = -1707.8334
$ws.Range("H22").Value2 = 1367.7858
$ws.Range("I22").Value2 = 1228.6666
$ws.Range("J22").Value2 = 1405.7273
$ws.Range("K22").Value2 = 1228.6666
$ws.Range("L22").Value2 = 1405.7273
$ws.Range("M22").Value2 = -878.6666
$ws.Range("N22").Value2 = -2105.7273
$ws.Range("H50").Value2 = 33999
$ws.Range("J50").Value2 = 48498.5
$ws.Range("L50").Value2 = 48498.5
$ws.Range("N50").Value2 = -49748.5
$ws.Range("H99").Value2 = 205487.88
$ws.Range("I99").Value2 = 504352.2
$ws.Range("J99").Value2 = 6245
$ws.Range("K99").Value2 = 504352.2
$ws.Range("L99").Value2 = 6245
$ws.Range("M99").Value2 = -502854.2
$ws.Range("N99").Value2 = -9241
$ws.Range("H107").Value2 = 10553.782
$ws.Range("I107").Value2 = 10965.318
$ws.Range("J107").Value2 = 1500
$ws.Range("K107").Value2 = 10965.318
$ws.Range("L107").Value2 = 1500
$ws.Range("M107").Value2 = -9045.317999999999
$ws.Range("N107").Value2 = -5340
$ws.Range("H113").Value2 = 1995.2858
$ws.Range("I113").Value2 = 1994.8334
$ws.Range("K113").Value2 = 1994.8334
$ws.Range("M113").Value2 = 175.1666
$ws.Range("H122").Value2 = 10144.643
$ws.Range("I122").Value2 = 18449.857
$ws.Range("J122").Value2 = 1839.4286
$ws.Range("K122").Value2 = 55349.571
$ws.Range("L122").Value2 = 5518.2858
$ws.Range("M122").Value2 = -52899.571
$ws.Range("N122").Value2 = -10418.2858
$ws.Range("H126").Value2 = 205487.88
$ws.Range("I126").Value2 = 504352.2
$ws.Range("J126").Value2 = 6245
$ws.Range("K126").Value2 = 1513056.6
$ws.Range("L126").Value2 = 18735
$ws.Range("M126").Value2 = -1510586.6
$ws.Range("N126").Value2 = -23675

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 26652196
$ws.Range("I4").Value2 = 4463239
$ws.Range("K4").Value2 = 13389717
$ws.Range("M4").Value2 = -13389605
$ws.Range("H14").Value2 = 537.25
$ws.Range("I14").Value2 = 537.25
$ws.Range("K14").Value2 = 1611.75
$ws.Range("M14").Value2 = -1438.75
$ws.Range("H114").Value2 = 4874.5
$ws.Range("J114").Value2 = 5666.3335
$ws.Range("L114").Value2 = 16999.0005
$ws.Range("N114").Value2 = -23507.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 6025.4
$ws.Range("I70").Value2 = 5225.4707
$ws.Range("J70").Value2 = 7725.25
$ws.Range("K70").Value2 = 5225.4707
$ws.Range("L70").Value2 = 7725.25
$ws.Range("M70").Value2 = -4955.4707
$ws.Range("N70").Value2 = -8265.25
$ws.Range("H73").Value2 = 6025.4
$ws.Range("I73").Value2 = 5225.4707
$ws.Range("J73").Value2 = 7725.25
$ws.Range("K73").Value2 = 5225.4707
$ws.Range("L73").Value2 = 7725.25
$ws.Range("M73").Value2 = -4289.4707
$ws.Range("N73").Value2 = -9597.25
$ws.Range("H107").Value2 = 340.7143
$ws.Range("I107").Value2 = 443
$ws.Range("J107").Value2 = 85
$ws.Range("K107").Value2 = 443
$ws.Range("L107").Value2 = 85
$ws.Range("M107").Value2 = 1477
$ws.Range("N107").Value2 = -3925
$ws.Range("H122").Value2 = 6381.41
$ws.Range("I122").Value2 = 3790.6177
$ws.Range("K122").Value2 = 11371.8531
$ws.Range("M122").Value2 = -8921.8531

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 14772.267
$ws.Range("I22").Value2 = 33900.168
$ws.Range("K22").Value2 = 33900.168
$ws.Range("M22").Value2 = -33605.168
$ws.Range("H27").Value2 = 14772.267
$ws.Range("I27").Value2 = 33900.168
$ws.Range("K27").Value2 = 33900.168
$ws.Range("M27").Value2 = -33793.168
$ws.Range("H100").Value2 = 4513.5
$ws.Range("I100").Value2 = 2866.5
$ws.Range("K100").Value2 = 2866.5
$ws.Range("M100").Value2 = -2325.5
$ws.Range("H132").Value2 = 682231.0600000001
$ws.Range("I132").Value2 = 1245341.4
$ws.Range("J132").Value2 = 6498.7
$ws.Range("K132").Value2 = 3736024.2
$ws.Range("L132").Value2 = 19496.1
$ws.Range("M132").Value2 = -3733494.2
$ws.Range("N132").Value2 = -24556.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value2 = 59989
$ws.Range("J133").Value2 = 59989
$ws.Range("L133").Value2 = 59989
$ws.Range("N133").Value2 = -70109
